# daily auto push: 2026-01-24 02:24 UTC
# Insert a new data row right after the existing "2026/01/24" (row 686),
# pushing rows 686..727 down to 687..728, then populate the newly
# inserted row 686 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 686 (and everything below it) down by one row.
$ws.Rows.Item(686).Insert()

# Force column A's text formatting on this single cell so the
# "yyyy/mm/dd"-shaped string is stored as literal text rather than being
# auto-parsed into a date serial by the smart Value setter, then restore
# the default "Normal" style so no stray number-format/style survives on
# the cell (matches the rest of the data rows, which carry no explicit
# style attribute).
$ws.Cells.Item(686, 1).NumberFormat = "@"
$ws.Cells.Item(686, 1).Value = "2026/01/24"
$ws.Cells.Item(686, 1).Style = "Normal"

$ws.Cells.Item(686, 2).Value = "土"
$ws.Cells.Item(686, 3).Value = 7
$ws.Cells.Item(686, 4).Value = 18
